$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05656054044545522
$ws.Range("D2").Value = 0.2075482465267271
$ws.Range("E2").Value = 0.2552946381726144
$ws.Range("F2").Value = 3.333329663450471
$ws.Range("G2").Value = 0.002583975018539875
$ws.Range("I2").Value = 2.395642200967131
$ws.Range("J2").Value = 0.3588304544780385
$ws.Range("M2").Value = 11.41472593379063

$ws.Range("C3").Value = 0.06090809120972329
$ws.Range("D3").Value = 0.1806669943574093
$ws.Range("E3").Value = 0.2237598830503629
$ws.Range("F3").Value = 3.445213852228918
$ws.Range("G3").Value = 0.002599019003974391
$ws.Range("I3").Value = 2.440428772863896
$ws.Range("J3").Value = 0.3186299615547057
$ws.Range("M3").Value = 10.15824920377452

$ws.Range("C4").Value = 0.06374839541784327
$ws.Range("D4").Value = 0.1642245098005333
$ws.Range("E4").Value = 0.2044081188946194
$ws.Range("F4").Value = 3.521793070907307
$ws.Range("G4").Value = 0.002608633122860526
$ws.Range("I4").Value = 2.473404507596712
$ws.Range("J4").Value = 0.2941527734772649
$ws.Range("M4").Value = 9.386271719037552

$ws.Range("C5").Value = 0.06494785784233592
$ws.Range("D5").Value = 0.1575373916587353
$ws.Range("E5").Value = 0.1965229795758745
$ws.Range("F5").Value = 3.554935876329495
$ws.Range("G5").Value = 0.002612646689714803
$ws.Range("I5").Value = 2.488184804421849
$ws.Range("J5").Value = 0.2842253406145403
$ws.Range("M5").Value = 9.071478596825898

$ws.Range("C6").Value = 0.06514953797769207
$ws.Range("D6").Value = 0.1564277365361875
$ws.Range("E6").Value = 0.1952136598363268
$ws.Range("F6").Value = 3.560554839942441
$ws.Range("G6").Value = 0.002613318947923631
$ws.Range("I6").Value = 2.490719145431868
$ws.Range("J6").Value = 0.2825796220681127
$ws.Range("M6").Value = 9.019192743556459

$ws.Range("C7").Value = 0.06376440290559549
$ws.Range("D7").Value = 0.1641342741693563
$ws.Range("E7").Value = 0.2043017759620938
$ws.Range("F7").Value = 3.522232266429882
$ws.Range("G7").Value = 0.002608686862348194
$ws.Range("I7").Value = 2.473598449516928
$ws.Range("J7").Value = 0.2940187032817221
$ws.Range("M7").Value = 9.382027245894221

$ws.Range("C8").Value = 0.05802352198778715
$ws.Range("D8").Value = 0.1982650624302664
$ws.Range("E8").Value = 0.2444180730842334
$ws.Range("F8").Value = 3.370243449847294
$ws.Range("G8").Value = 0.002589084494513873
$ws.Range("I8").Value = 2.4099267108395
$ws.Range("J8").Value = 0.3449236650505441
$ws.Range("M8").Value = 10.98154423438388

$ws.Range("C9").Value = 0.04816613556701599
$ws.Range("D9").Value = 0.2658117973180083
$ws.Range("E9").Value = 0.3232654723109789
$ws.Range("F9").Value = 3.136862360837881
$ws.Range("G9").Value = 0.002553592284770111
$ws.Range("I9").Value = 2.330138991276542
$ws.Range("J9").Value = 0.4466134486579847
$ws.Range("M9").Value = 14.11841718582878

$ws.Range("C10").Value = 0.04184178325499843
$ws.Range("D10").Value = 0.3159919893975882
$ws.Range("E10").Value = 0.3814533064968373
$ws.Range("F10").Value = 3.007879980196037
$ws.Range("G10").Value = 0.002529249383284637
$ws.Range("I10").Value = 2.301326257851684
$ws.Range("J10").Value = 0.5228161314431645
$ws.Range("M10").Value = 16.42971922713639

$ws.Range("C11").Value = 0.03917890693996462
$ws.Range("D11").Value = 0.3389823449212201
$ws.Range("E11").Value = 0.4080174662851448
$ws.Range("F11").Value = 2.959133256360872
$ws.Range("G11").Value = 0.002518536807825367
$ws.Range("I11").Value = 2.29522883347579
$ws.Range("J11").Value = 0.5578929729263393
$ws.Range("M11").Value = 17.48422909488204

$ws.Range("C12").Value = 0.03820266581377396
$ws.Range("D12").Value = 0.3477152502007073
$ws.Range("E12").Value = 0.4180934067476159
$ws.Range("F12").Value = 2.942163918013961
$ws.Range("G12").Value = 0.00251453091751811
$ws.Range("I12").Value = 2.293975057569227
$ws.Range("J12").Value = 0.5712425848237217
$ws.Range("M12").Value = 17.88412875839322

$ws.Range("C13").Value = 0.03841146546659413
$ws.Range("D13").Value = 0.3458332059236398
$ws.Range("E13").Value = 0.4159225787832526
$ws.Range("F13").Value = 2.945751273035881
$ws.Range("G13").Value = 0.002515391421163126
$ws.Range("I13").Value = 2.294197372680628
$ws.Range("J13").Value = 0.5683644106220527
$ws.Range("M13").Value = 17.79797514926526

$ws.Range("C14").Value = 0.03909793903581615
$ws.Range("D14").Value = 0.3397002452521463
$ws.Range("E14").Value = 0.4088460670711953
$ws.Range("F14").Value = 2.95770694826146
$ws.Range("G14").Value = 0.002518206231185914
$ws.Range("I14").Value = 2.295104257829991
$ws.Range("J14").Value = 0.5589898699062701
$ws.Range("M14").Value = 17.51711664633206

$ws.Range("C15").Value = 0.03952264973267638
$ws.Range("D15").Value = 0.3359472477254712
$ws.Range("E15").Value = 0.4045137673118262
$ws.Range("F15").Value = 2.965226091377076
$ws.Range("G15").Value = 0.002519936953525826
$ws.Range("I15").Value = 2.295798610881846
$ws.Range("J15").Value = 0.5532566232052147
$ws.Range("M15").Value = 17.34516241147986

$ws.Range("C16").Value = 0.04202022997396249
$ws.Range("D16").Value = 0.3144930925328708
$ws.Range("E16").Value = 0.3797194329811475
$ws.Range("F16").Value = 3.011270726707664
$ws.Range("G16").Value = 0.002529956641052987
$ws.Range("I16").Value = 2.301869770705366
$ws.Range("J16").Value = 0.5205326966751613
$ws.Range("M16").Value = 16.36087799942749

$ws.Range("C17").Value = 0.04360818417469314
$ws.Range("D17").Value = 0.3013758691899113
$ws.Range("E17").Value = 0.364535182913869
$ws.Range("F17").Value = 3.042103776365423
$ws.Range("G17").Value = 0.002536195053266764
$ws.Range("I17").Value = 2.307422394272635
$ws.Range("J17").Value = 0.5005682982596511
$ws.Range("M17").Value = 15.75793518329931

$ws.Range("C18").Value = 0.04454160093062587
$ws.Range("D18").Value = 0.2938463225937085
$ws.Range("E18").Value = 0.3558102701956614
$ws.Range("F18").Value = 3.060768684445122
$ws.Range("G18").Value = 0.002539817289642099
$ws.Range("I18").Value = 2.311272930239568
$ws.Range("J18").Value = 0.4891235833190706
$ws.Range("M18").Value = 15.41142020399315

$ws.Range("C19").Value = 0.04486104669500079
$ws.Range("D19").Value = 0.2912994434676079
$ws.Range("E19").Value = 0.3528575622702022
$ws.Range("F19").Value = 3.067246277433668
$ws.Range("G19").Value = 0.002541049603151943
$ws.Range("I19").Value = 2.312688064944297
$ws.Range("J19").Value = 0.4852549536334436
$ws.Range("M19").Value = 15.2941406244264

$ws.Range("C20").Value = 0.04343705733557179
$ws.Range("D20").Value = 0.3027706290929189
$ws.Range("E20").Value = 0.366150653550136
$ws.Range("F20").Value = 3.038724795483006
$ws.Range("G20").Value = 0.002535527446218877
$ws.Range("I20").Value = 2.306762996088366
$ws.Range("J20").Value = 0.5026895280438737
$ws.Range("M20").Value = 15.8220893636817

$ws.Range("C21").Value = 0.03889542194132112
$ws.Range("D21").Value = 0.341500883940796
$ws.Range("E21").Value = 0.4109241315067038
$ws.Range("F21").Value = 2.954154325573882
$ws.Range("G21").Value = 0.002517378086256727
$ws.Range("I21").Value = 2.294808857678618
$ws.Range("J21").Value = 0.5617415243887365
$ws.Range("M21").Value = 17.59959472766877

$ws.Range("C22").Value = 0.03611517376221407
$ws.Range("D22").Value = 0.3669722457311764
$ws.Range("E22").Value = 0.4402849908347974
$ws.Range("F22").Value = 2.907598764815589
$ws.Range("G22").Value = 0.002505811519482462
$ws.Range("I22").Value = 2.293170293705515
$ws.Range("J22").Value = 0.6007285657380805
$ws.Range("M22").Value = 18.76474702840528

$ws.Range("C23").Value = 0.03758138005613532
$ws.Range("D23").Value = 0.3533619494351399
$ws.Range("E23").Value = 0.4246044066457557
$ws.Range("F23").Value = 2.931627347223753
$ws.Range("G23").Value = 0.002511958225769427
$ws.Range("I23").Value = 2.293463636507624
$ws.Range("J23").Value = 0.5798818217844541
$ws.Range("M23").Value = 18.14251987880368

$ws.Range("C24").Value = 0.04351436008688836
$ws.Range("D24").Value = 0.302140021898424
$ws.Range("E24").Value = 0.3654202848224202
$ws.Range("F24").Value = 3.040249513960021
$ws.Range("G24").Value = 0.002535829159915944
$ws.Range("I24").Value = 2.307059063220549
$ws.Range("J24").Value = 0.5017304175736115
$ws.Range("M24").Value = 15.79308487989636

$ws.Range("C25").Value = 0.0506764719789059
$ws.Range("D25").Value = 0.2474548565953114
$ws.Range("E25").Value = 0.3019013171180518
$ws.Range("F25").Value = 3.192767414498519
$ws.Range("G25").Value = 0.002562884707639039
$ws.Range("I25").Value = 2.346664137761238
$ws.Range("J25").Value = 0.4188691494290424
$ws.Range("M25").Value = 13.26914682592604
